$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95
$ws.Range("A95").Value = 130964396
$ws.Range("B95").Value = 98930
$ws.Range("E95").Value = 219790
$ws.Range("Q95").Value = 509365
$ws.Range("R95").Value = 6718883
$ws.Range("D95").Value = "LC"
$ws.Range("F95").Value = "Fläcknycklar"
$ws.Range("G95").Value = "Dactylorhiza maculata"
$ws.Range("H95").Value = "(L.) Soó"
$ws.Range("Y95").NumberFormat = "@"
$ws.Range("Y95").Value = "2025-07-03"
$ws.Range("AA95").NumberFormat = "@"
$ws.Range("AA95").Value = "2025-07-03"
$ws.Range("AC95").Value = "Betydande förekomster . inventering åt vasa vind"

# Row 96
$ws.Range("A96").Value = 130964535
$ws.Range("Q96").Value = 509939
$ws.Range("R96").Value = 6719007
$ws.Range("AC96").Value = "Flera . inventering åt vasa vind"

# Row 97
$ws.Range("A97").Value = 130964529
$ws.Range("Q97").Value = 509588
$ws.Range("R97").Value = 6719099
$ws.Range("AC97").Value = "Enstaka . inventering åt vasa vind"

# Row 98
$ws.Range("A98").Value = 130964573
$ws.Range("B98").Value = 79243
$ws.Range("E98").Value = 6425
$ws.Range("Q98").Value = 509515
$ws.Range("R98").Value = 6719063
$ws.Range("D98").Value = "NT"
$ws.Range("F98").Value = "Garnlav"
$ws.Range("G98").Value = "Alectoria sarmentosa"
$ws.Range("H98").Value = "(Ach.) Ach."
$ws.Range("Y98").NumberFormat = "@"
$ws.Range("Y98").Value = "2025-07-02"
$ws.Range("AA98").NumberFormat = "@"
$ws.Range("AA98").Value = "2025-07-02"
$ws.Range("AC98").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 104
$ws.Range("A104").Value = 130964545
$ws.Range("B104").Value = 57073
$ws.Range("E104").Value = 100138
$ws.Range("Q104").Value = 509535
$ws.Range("R104").Value = 6718925
$ws.Range("D104").Value = "LC"
$ws.Range("F104").Value = "Tjäder"
$ws.Range("G104").Value = "Tetrao urogallus"
$ws.Range("H104").Value = "Linnaeus, 1758"
$ws.Range("AC104").Value = "Spillning . inventering åt vasa vind"

# Row 105
$ws.Range("A105").Value = 130964541
$ws.Range("B105").Value = 91808
$ws.Range("E105").Value = 1202
$ws.Range("Q105").Value = 509703
$ws.Range("R105").Value = 6719018
$ws.Range("D105").Value = "NT"
$ws.Range("F105").Value = "Ullticka"
$ws.Range("G105").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H105").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("AC105").Value = "Enstaka . inventering åt vasa vind"

# Row 109
$ws.Range("A109").Value = 130964647
$ws.Range("B109").Value = 92106
$ws.Range("E109").Value = 658
$ws.Range("Q109").Value = 509741
$ws.Range("R109").Value = 6718998
$ws.Range("D109").Value = "NT"
$ws.Range("F109").Value = "Rosenticka"
$ws.Range("G109").Value = "Fomitopsis rosea"
$ws.Range("H109").Value = "(Alb. & Schwein.:Fr.) P.Karst."
$ws.Range("AC109").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 110
$ws.Range("A110").Value = 130964641
$ws.Range("B110").Value = 98930
$ws.Range("E110").Value = 219790
$ws.Range("Q110").Value = 509932
$ws.Range("R110").Value = 6719045
$ws.Range("D110").Value = "LC"
$ws.Range("F110").Value = "Fläcknycklar"
$ws.Range("G110").Value = "Dactylorhiza maculata"
$ws.Range("H110").Value = "(L.) Soó"
$ws.Range("AC110").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 115
$ws.Range("A115").Value = 130964539
$ws.Range("B115").Value = 98930
$ws.Range("E115").Value = 219790
$ws.Range("Q115").Value = 509808
$ws.Range("R115").Value = 6719041
$ws.Range("D115").Value = "LC"
$ws.Range("F115").Value = "Fläcknycklar"
$ws.Range("G115").Value = "Dactylorhiza maculata"
$ws.Range("H115").Value = "(L.) Soó"
$ws.Range("AC115").Value = "Flera . inventering åt vasa vind"

# Row 116
$ws.Range("A116").Value = 130964527
$ws.Range("B116").Value = 79243
$ws.Range("E116").Value = 6425
$ws.Range("Q116").Value = 509597
$ws.Range("R116").Value = 6719076
$ws.Range("D116").Value = "NT"
$ws.Range("F116").Value = "Garnlav"
$ws.Range("G116").Value = "Alectoria sarmentosa"
$ws.Range("H116").Value = "(Ach.) Ach."
$ws.Range("AC116").Value = "Rikligt . inventering åt vasa vind"

# Row 119
$ws.Range("A119").Value = 130964649
$ws.Range("B119").Value = 98930
$ws.Range("E119").Value = 219790
$ws.Range("Q119").Value = 509705
$ws.Range("R119").Value = 6718923
$ws.Range("D119").Value = "LC"
$ws.Range("F119").Value = "Fläcknycklar"
$ws.Range("G119").Value = "Dactylorhiza maculata"
$ws.Range("H119").Value = "(L.) Soó"
$ws.Range("AC119").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 120
$ws.Range("A120").Value = 130964648
$ws.Range("B120").Value = 92267
$ws.Range("E120").Value = 1209
$ws.Range("Q120").Value = 509744
$ws.Range("R120").Value = 6718982
$ws.Range("D120").Value = "VU"
$ws.Range("F120").Value = "Rynkskinn"
$ws.Range("G120").Value = "Hermanssonia centrifuga"
$ws.Range("H120").Value = "(P. Karst.) Zmitr."
$ws.Range("AC120").Value = "Måttliga förekomster . inventering åt vasa vind"

# Row 122
$ws.Range("A122").Value = 130964644
$ws.Range("B122").Value = 98917
$ws.Range("E122").Value = 220093
$ws.Range("Q122").Value = 509801
$ws.Range("R122").Value = 6719017
$ws.Range("F122").Value = "Korallrot"
$ws.Range("G122").Value = "Corallorhiza trifida"
$ws.Range("H122").Value = "Châtel."
$ws.Range("AC122").Value = "Sparsamma förekomster . inventering åt vasa vind"

# Row 123
$ws.Range("A123").Value = 130964542
$ws.Range("B123").Value = 57073
$ws.Range("E123").Value = 100138
$ws.Range("Q123").Value = 509635
$ws.Range("R123").Value = 6718941
$ws.Range("D123").Value = "LC"
$ws.Range("F123").Value = "Tjäder"
$ws.Range("G123").Value = "Tetrao urogallus"
$ws.Range("H123").Value = "Linnaeus, 1758"
$ws.Range("AC123").Value = "Spillning . inventering åt vasa vind"

# Row 124
$ws.Range("A124").Value = 130964640
$ws.Range("B124").Value = 57881
$ws.Range("E124").Value = 100049
$ws.Range("Q124").Value = 509697
$ws.Range("R124").Value = 6719144
$ws.Range("D124").Value = "NT"
$ws.Range("F124").Value = "Spillkråka"
$ws.Range("G124").Value = "Dryocopus martius"
$ws.Range("H124").Value = "(Linnaeus, 1758)"
$ws.Range("AC124").Value = "Gamla födosöksspår . inventering åt vasa vind"

# Row 125
$ws.Range("A125").Value = 130964643
$ws.Range("B125").Value = 98930
$ws.Range("E125").Value = 219790
$ws.Range("Q125").Value = 509829
$ws.Range("R125").Value = 6719000
$ws.Range("D125").Value = "LC"
$ws.Range("F125").Value = "Fläcknycklar"
$ws.Range("G125").Value = "Dactylorhiza maculata"
$ws.Range("H125").Value = "(L.) Soó"
$ws.Range("Y125").NumberFormat = "@"
$ws.Range("Y125").Value = "2025-07-02"
$ws.Range("AA125").NumberFormat = "@"
$ws.Range("AA125").Value = "2025-07-02"
$ws.Range("AC125").Value = "Måttlig förekomst . inventering åt vasa vind"

# Row 126
$ws.Range("A126").Value = 130964390
$ws.Range("B126").Value = 99013
$ws.Range("E126").Value = 220787
$ws.Range("Q126").Value = 509475
$ws.Range("R126").Value = 6718881
$ws.Range("D126").Value = "VU"
$ws.Range("F126").Value = "Knärot"
$ws.Range("G126").Value = "Goodyera repens"
$ws.Range("H126").Value = "(L.) R. Br."
$ws.Range("Y126").NumberFormat = "@"
$ws.Range("Y126").Value = "2025-07-03"
$ws.Range("AA126").NumberFormat = "@"
$ws.Range("AA126").Value = "2025-07-03"
$ws.Range("AC126").Value = "Måttliga förekomster, Ca 10-15 plantor . inventering åt vasa vind"

# Row 127
$ws.Range("A127").Value = 130964538
$ws.Range("B127").Value = 79243
$ws.Range("E127").Value = 6425
$ws.Range("Q127").Value = 509875
$ws.Range("R127").Value = 6719025
$ws.Range("D127").Value = "NT"
$ws.Range("F127").Value = "Garnlav"
$ws.Range("G127").Value = "Alectoria sarmentosa"
$ws.Range("H127").Value = "(Ach.) Ach."
$ws.Range("AC127").Value = "Enstaka . inventering åt vasa vind"

# Row 128
$ws.Range("A128").Value = 130964546
$ws.Range("B128").Value = 92503
$ws.Range("E128").Value = 898
$ws.Range("Q128").Value = 509515
$ws.Range("R128").Value = 6718886
$ws.Range("D128").Value = "VU"
$ws.Range("F128").Value = "Blackticka"
$ws.Range("G128").Value = "Steccherinum collabens"
$ws.Range("H128").Value = "(Fr.) Vesterholt"
$ws.Range("AC128").Value = "Betydande förekomst . inventering åt vasa vind"

# Row 129
$ws.Range("A129").Value = 130964639
$ws.Range("B129").Value = 57073
$ws.Range("E129").Value = 100138
$ws.Range("Q129").Value = 509645
$ws.Range("R129").Value = 6719169
$ws.Range("D129").Value = "LC"
$ws.Range("F129").Value = "Tjäder"
$ws.Range("G129").Value = "Tetrao urogallus"
$ws.Range("H129").Value = "Linnaeus, 1758"
$ws.Range("AC129").Value = "Vinterspillning . inventering åt vasa vind"

# Row 130
$ws.Range("A130").Value = 130964646
$ws.Range("B130").Value = 91808
$ws.Range("E130").Value = 1202
$ws.Range("Q130").Value = 509764
$ws.Range("R130").Value = 6719043
$ws.Range("D130").Value = "NT"
$ws.Range("F130").Value = "Ullticka"
$ws.Range("G130").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H130").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("AC130").Value = "Måttlig förekomst . inventering åt vasa vind"
